# repull data, push all data, mean calculation
# Update the dSF column (F) values for specific rows to reflect the
# re-pulled / recalculated data. Column E (dS0) stays as-is; only the
# F-column values listed below change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    8  = -8
    12 = 4
    13 = -2
    17 = -2
    18 = -5
    19 = -3
    20 = 1
    22 = -3
    23 = -10
    24 = -3
    25 = -1
    26 = -5
    28 = -3
    30 = 2
    31 = -6
    33 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
